# Apply "Finished initial version of fab files for 12-12 teensy controller
# arena" edits to the CPL sheet: re-measured placement/rotation data for the
# outer ring of parts (rows 58-69, which also grow two new formatted-but-
# empty columns F/G), a handful of rotation-only fixes (rows 38, 53-57), and
# the resulting selection state left behind in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rotation-only corrections -------------------------------------------------
$ws.Range("E38").Value = 0

$ws.Range("E53").Value = -90
$ws.Range("E54").Value = -90
$ws.Range("E55").Value = 90
$ws.Range("E56").Value = 90
$ws.Range("E57").Value = 90

# --- Re-measured placements + rotations for the outer ring (rows 58-69) -------
$ws.Range("B58").Value = 278.659028
$ws.Range("E58").Value = 180

$ws.Range("B59").Value = 274.247981
$ws.Range("C59").Value = -143.134745
$ws.Range("E59").Value = 210

$ws.Range("B60").Value = 243.993774
$ws.Range("C60").Value = -99.555018
$ws.Range("E60").Value = 240

$ws.Range("C61").Value = -76.940971
$ws.Range("E61").Value = 270

$ws.Range("B62").Value = 143.134745
$ws.Range("C62").Value = -81.352018
$ws.Range("E62").Value = 300

$ws.Range("B63").Value = 99.555018
$ws.Range("C63").Value = -111.606225
$ws.Range("E63").Value = 330

$ws.Range("B64").Value = 76.940971
$ws.Range("E64").Value = 0

$ws.Range("B65").Value = 81.352018
$ws.Range("C65").Value = -212.465254
$ws.Range("E65").Value = 30

$ws.Range("B66").Value = 111.606225
$ws.Range("C66").Value = -256.044981
$ws.Range("E66").Value = 60

$ws.Range("C67").Value = -278.659028
$ws.Range("E67").Value = 90

$ws.Range("B68").Value = 212.465254
$ws.Range("C68").Value = -274.247981
$ws.Range("E68").Value = 120

$ws.Range("B69").Value = 256.044981
$ws.Range("C69").Value = -243.993774
$ws.Range("E69").Value = 150

# New (blank) F & G columns for rows 58-69, carrying the same cell format as
# column E on those rows (picked up via a format-only paste so the existing
# style is reused rather than a fresh one minted).
$ws.Range("E58:E69").Copy()
$ws.Range("F58:F69").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E58:E69").Copy()
$ws.Range("G58:G69").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Leftover selection state from the editing session -------------------------
$ws.Range("H14").Select()
